$d = $word.ActiveDocument

# 1) Update the "Data de envio" timestamp line.
$d.Content.Find.Execute("Data de envio: 23/06/2025, 14:46:53", $true, $false, $false, $false, $false, $true, 1, $false, "Data de envio: 23/06/2025, 19:17:33", 2) | Out-Null

# 2) Update each labeled field value by locating the paragraph whose
#    bold label run matches, then replacing the trailing value text
#    that follows the line break within that same paragraph.
$fieldUpdates = @(
    @{ Label = "Nome completo da empresa:"; OldValue = " Não informado"; NewValue = " dasdsa" }
    @{ Label = "Outros nomes/apelidos:"; OldValue = " Não informado"; NewValue = " dsadsa" }
    @{ Label = "Definição do negócio:"; OldValue = " Não informado"; NewValue = " dsa" }
    @{ Label = "CNPJ:"; OldValue = " Não informado"; NewValue = " dasdsa" }
    @{ Label = "Categoria do negócio:"; OldValue = " Não informado"; NewValue = " dasdsa" }
    @{ Label = "Endereço completo:"; OldValue = " Não informado"; NewValue = " dsadsadas" }
    @{ Label = "Tipo de estabelecimento:"; OldValue = " Não informado"; NewValue = " Online/Virtual" }
    @{ Label = "Identificação na fachada:"; OldValue = " Não informado"; NewValue = " Sim, bem visível" }
    @{ Label = "Áreas de atendimento:"; OldValue = " Não informado"; NewValue = " assa" }
    @{ Label = "Raio de atendimento:"; OldValue = " Não informado"; NewValue = " saddas" }
    @{ Label = "Taxa de deslocamento:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Horários de funcionamento:"; OldValue = " Não informado"; NewValue = " dasdsadaw" }
    @{ Label = "Horários especiais:"; OldValue = " Não informado"; NewValue = " wdwaedwa" }
    @{ Label = "História da criação:"; OldValue = " Não informado"; NewValue = " dasdas" }
    @{ Label = "Experiência no ramo:"; OldValue = " Não informado"; NewValue = " dasdasdasd" }
    @{ Label = "Certificações:"; OldValue = " Não informado"; NewValue = " asdasdas" }
    @{ Label = "Prêmios e reconhecimentos:"; OldValue = " Não informado"; NewValue = " dasdas" }
    @{ Label = "Quantidade de clientes:"; OldValue = " Não informado"; NewValue = " Mais de 1000" }
    @{ Label = "Lista de produtos/serviços:"; OldValue = " Não informado"; NewValue = " dasdassda" }
    @{ Label = "Carro-chefe:"; OldValue = " Não informado"; NewValue = " asddas" }
    @{ Label = "Produtos sazonais:"; OldValue = " Não informado"; NewValue = " dsadasdasd" }
    @{ Label = "Diferencial da concorrência:"; OldValue = " Não informado"; NewValue = " asdasdas" }
    @{ Label = "Marcas parceiras:"; OldValue = " Não informado"; NewValue = " dasdsa" }
    @{ Label = "Faixa de preço:"; OldValue = " Não informado"; NewValue = " Econômico/Popular" }
    @{ Label = "Garantia:"; OldValue = " Não informado"; NewValue = " dasdsadas" }
    @{ Label = "Canais de compra:"; OldValue = " Não informado"; NewValue = " dasdasdsa" }
    @{ Label = "Contato preferencial:"; OldValue = " Não informado"; NewValue = " dsadasdas" }
    @{ Label = "Agendamento online:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Formas de pagamento:"; OldValue = " Não informado"; NewValue = " dasdsa" }
    @{ Label = "Valor mínimo:"; OldValue = " Não informado"; NewValue = " dasdas" }
    @{ Label = "Tempo de entrega:"; OldValue = " Não informado"; NewValue = " dassadas" }
    @{ Label = "Instruções especiais:"; OldValue = " Não informado"; NewValue = " dasdasdsadas" }
    @{ Label = "Como clientes procuram:"; OldValue = " Não informado"; NewValue = " dasdwedaew" }
    @{ Label = "Palavras-chave desejadas:"; OldValue = " Não informado"; NewValue = " aeadasdasdas" }
    @{ Label = "Termos a evitar:"; OldValue = " Não informado"; NewValue = " dasdsadsa" }
    @{ Label = "Busca principal:"; OldValue = " Não informado"; NewValue = " Nome da empresa" }
    @{ Label = "Concorrentes diretos:"; OldValue = " Não informado"; NewValue = " dsadsadsadsa" }
    @{ Label = "O que admira nos concorrentes:"; OldValue = " Não informado"; NewValue = " dasdsadas" }
    @{ Label = "O que faz melhor:"; OldValue = " Não informado"; NewValue = " dassdadsa" }
    @{ Label = "Melhor presença digital:"; OldValue = " Não informado"; NewValue = " dasdsadsa" }
    @{ Label = "Benchmark externo:"; OldValue = " Não informado"; NewValue = " dsadasdsa" }
    @{ Label = "Identidade visual:"; OldValue = " Manual de marca completo"; NewValue = " Logo profissional completo" }
    @{ Label = "Tipos de fotos disponíveis:"; OldValue = " Não informado"; NewValue = " dsadsa" }
    @{ Label = "Possui vídeos:"; OldValue = " Pretende fazer"; NewValue = " Sim" }
    @{ Label = "Autorização para fotos:"; OldValue = " Com autorização específica"; NewValue = " Sim" }
    @{ Label = "Redes sociais:"; OldValue = " Não informado"; NewValue = " sdsdasda" }
    @{ Label = "Site próprio:"; OldValue = " Não informado"; NewValue = " saddsa" }
    @{ Label = "Plataformas presentes:"; OldValue = " Não informado"; NewValue = " dasdsadsa" }
    @{ Label = "Google Ads:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Integrar plataformas no GMB:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Atributos do negócio:"; OldValue = " Não informado"; NewValue = " dsadsadsadas" }
    @{ Label = "Características do ambiente:"; OldValue = " Não informado"; NewValue = " dasdsadsa" }
    @{ Label = "Público-alvo:"; OldValue = " Não informado"; NewValue = " Famílias" }
    @{ Label = "Avaliações online recebidas:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Onde recebeu avaliações:"; OldValue = " Não informado"; NewValue = " dsadsadsa" }
    @{ Label = "Estratégia para avaliações negativas:"; OldValue = " Não informado"; NewValue = " dasdasdas" }
    @{ Label = "Estratégia para solicitar avaliações:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Objetivo principal:"; OldValue = " Não informado"; NewValue = " Gerar mais leads" }
    @{ Label = "Meta de clientes mensais:"; OldValue = " Não informado"; NewValue = " 11-30" }
    @{ Label = "Google Ads futuro:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Responsável pela gestão:"; OldValue = " Não informado"; NewValue = " dasdsadsa" }
    @{ Label = "Acesso ao email Google:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Tentativa anterior GMB:"; OldValue = " Não informado"; NewValue = " dsadsadsa" }
    @{ Label = "Autoriza cartão postal:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Informações a ocultar:"; OldValue = " Não informado"; NewValue = " dasdsadsa" }
    @{ Label = "Restrições legais:"; OldValue = " Não informado"; NewValue = " dsasda" }
    @{ Label = "Produtos restritos pelo Google:"; OldValue = " Não informado"; NewValue = " asddsadas" }
    @{ Label = "Problemas anteriores com Google:"; OldValue = " Não informado"; NewValue = " Sim" }
    @{ Label = "Detalhes importantes:"; OldValue = " Não informado"; NewValue = " dsadsadsa" }
    @{ Label = "Maior expectativa:"; OldValue = " Não informado"; NewValue = " dasdasdas" }
    @{ Label = "Orçamento para melhorias:"; OldValue = " Não informado"; NewValue = " Apenas configuração" }
)

foreach ($update in $fieldUpdates) {
    foreach ($para in $d.Paragraphs) {
        $pr = $para.Range
        if ($pr.Text.StartsWith($update.Label)) {
            $pr.Find.Execute($update.OldValue, $true, $false, $false, $false, $false, $true, 1, $false, $update.NewValue, 2) | Out-Null
            break
        }
    }
}
